# "test toegevoegd voor Join, Split en Paste"
# Adds test data/results for the Join, Split and Paste algorithms to the
# "Blad1" worksheet (columns H-L, rows 1-26), alongside the existing
# size/Find/Insert/Delete results in columns A-E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
# Write J1:L1 (Join/Split/Paste) before H1:I1 (left/right) so that new
# shared-string entries are appended in the same order used by the
# author: Join, Split, Paste, left, right.
$ws.Range("J1").Value = "Join"
$ws.Range("K1").Value = "Split"
$ws.Range("L1").Value = "Paste"
$ws.Range("H1").Value = "left"
$ws.Range("I1").Value = "right"

# --- Existing LOG()/*10 helper columns become "fill down" formulas --
# Re-entering them as one range fill lets Excel collapse the series
# into shared formulas, matching how the original columns were
# extended.
$ws.Range("A3:A9").Formula = "=+A2*10"
$ws.Range("E2:E9").Formula = "=LOG(A2)"

# --- New "left"/"right" driver columns (H:I) -------------------------
# H/I together enumerate every (left-size, right-size) combination used
# by the Join/Split/Paste benchmarks: right cycles 1,10,100,1000,10000
# and left advances by a factor of 10 every time right wraps back to 1.
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1

$ws.Range("H3").Formula = "=IF(I3=1,H2*10,H2)"
$ws.Range("I3").Formula = "=IF(I2=10000,1,I2*10)"

$ws.Range("H4:H16").Formula = "=IF(I4=1,H3*10,H3)"
$ws.Range("I4:I16").Formula = "=IF(I3=10000,1,I3*10)"

$ws.Range("H17:H24").Formula = "=IF(I17=1,H16*10,H16)"
$ws.Range("I17:I24").Formula = "=IF(I16=10000,1,I16*10)"

# Fill one row further (to 27) then remove it again: this reproduces
# the trailing shared-formula range (ref spans 25:27) while only rows
# 25 and 26 actually hold data, exactly as in the authored workbook.
$ws.Range("H25:H27").Formula = "=IF(I25=1,H24*10,H24)"
$ws.Range("I25:I27").Formula = "=IF(I24=10000,1,I24*10)"
$ws.Rows.Item(27).Delete()

# --- Measured timings for Join (J), Split (K) and Paste (L) ---------
$ws.Cells.Item(2, 10).Value = 251
$ws.Cells.Item(2, 11).Value = 529
$ws.Cells.Item(2, 12).Value = 301
$ws.Cells.Item(3, 10).Value = 230
$ws.Cells.Item(3, 11).Value = 643
$ws.Cells.Item(3, 12).Value = 198
$ws.Cells.Item(4, 10).Value = 129
$ws.Cells.Item(4, 11).Value = 352
$ws.Cells.Item(4, 12).Value = 117
$ws.Cells.Item(5, 10).Value = 182
$ws.Cells.Item(5, 11).Value = 358
$ws.Cells.Item(5, 12).Value = 127
$ws.Cells.Item(6, 10).Value = 564
$ws.Cells.Item(6, 11).Value = 442
$ws.Cells.Item(6, 12).Value = 156
$ws.Cells.Item(7, 10).Value = 62
$ws.Cells.Item(7, 11).Value = 157
$ws.Cells.Item(7, 12).Value = 84
$ws.Cells.Item(8, 10).Value = 74
$ws.Cells.Item(8, 11).Value = 272
$ws.Cells.Item(8, 12).Value = 109
$ws.Cells.Item(9, 10).Value = 109
$ws.Cells.Item(9, 11).Value = 327
$ws.Cells.Item(9, 12).Value = 127
$ws.Cells.Item(10, 10).Value = 139
$ws.Cells.Item(10, 11).Value = 364
$ws.Cells.Item(10, 12).Value = 150
$ws.Cells.Item(11, 10).Value = 579
$ws.Cells.Item(11, 11).Value = 528
$ws.Cells.Item(11, 12).Value = 231
$ws.Cells.Item(12, 10).Value = 77
$ws.Cells.Item(12, 11).Value = 242
$ws.Cells.Item(12, 12).Value = 115
$ws.Cells.Item(13, 10).Value = 102
$ws.Cells.Item(13, 11).Value = 316
$ws.Cells.Item(13, 12).Value = 154
$ws.Cells.Item(14, 10).Value = 124
$ws.Cells.Item(14, 11).Value = 389
$ws.Cells.Item(14, 12).Value = 198
$ws.Cells.Item(15, 10).Value = 190
$ws.Cells.Item(15, 11).Value = 446
$ws.Cells.Item(15, 12).Value = 245
$ws.Cells.Item(16, 10).Value = 678
$ws.Cells.Item(16, 11).Value = 606
$ws.Cells.Item(16, 12).Value = 246
$ws.Cells.Item(17, 10).Value = 104
$ws.Cells.Item(17, 11).Value = 365
$ws.Cells.Item(17, 12).Value = 144
$ws.Cells.Item(18, 10).Value = 133
$ws.Cells.Item(18, 11).Value = 369
$ws.Cells.Item(18, 12).Value = 200
$ws.Cells.Item(19, 10).Value = 186
$ws.Cells.Item(19, 11).Value = 460
$ws.Cells.Item(19, 12).Value = 218
$ws.Cells.Item(20, 10).Value = 218
$ws.Cells.Item(20, 11).Value = 568
$ws.Cells.Item(20, 12).Value = 251
$ws.Cells.Item(21, 10).Value = 818
$ws.Cells.Item(21, 11).Value = 628
$ws.Cells.Item(21, 12).Value = 338
$ws.Cells.Item(22, 10).Value = 170
$ws.Cells.Item(22, 11).Value = 432
$ws.Cells.Item(22, 12).Value = 178
$ws.Cells.Item(23, 10).Value = 204
$ws.Cells.Item(23, 11).Value = 489
$ws.Cells.Item(23, 12).Value = 255
$ws.Cells.Item(24, 10).Value = 246
$ws.Cells.Item(24, 11).Value = 596
$ws.Cells.Item(24, 12).Value = 263
$ws.Cells.Item(25, 10).Value = 380
$ws.Cells.Item(25, 11).Value = 672
$ws.Cells.Item(25, 12).Value = 300
$ws.Cells.Item(26, 10).Value = 817
$ws.Cells.Item(26, 11).Value = 789
$ws.Cells.Item(26, 12).Value = 391

# --- Selection / active sheet ----------------------------------------
# The author finished editing with cell O24 selected on Blad1, and the
# Blad1 tab (rather than the Grafiek1 chart tab) active.
$ws.Activate()
$ws.Range("O24").Select()
